$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(-0.299, 0.175, -0.383, 118, 69, -163),
    @(-0.091, -0.063, -0.054, -89, -63, -146),
    @(-0.13, -0.056, -0.297, -54, 44, 8),
    @(-0.769, -0.29, -0.073, 173, 43, -101),
    @(0.131, -0.145, 0.282, -62, -43, -140),
    @(-0.167, -0.129, 0.107, 85, -52, -12),
    @(-0.236, 0.13, 0.013, 111, 7, 75),
    @(0.218, -0.165, 0.914, -72, 36, -47),
    @(-0.579, 0.202, -0.078, -110, -7, 133),
    @(0.091, -0.057, 0.056, 8, 48, 116)
)

for ($r = 1; $r -le 10; $r++) {
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$r - 1][$c - 1]
    }
}
